$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "2023-12-06 08:28:48"
$ws.Range("B9").Value = 0.0008

$ws.Range("A10").Value = "2023-12-06 08:29:13"
$ws.Range("B10").Value = 0.001

$ws.Range("A11").Value = "2023-12-06 08:30:26"
$ws.Range("B11").Value = 0.0038
